$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 97, shifting existing rows 97:100 down to 98:101.
$ws.Rows("97:97").Insert()

# Populate the newly inserted row 97 with the new weekly price record.
$ws.Range("A97").Value = 1
$ws.Range("B97").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C97").Value = "Arica y Parinacota"
$ws.Range("D97").Value = "10/13/2023"
$ws.Range("E97").Value = 15
$ws.Range("F97").Value = 100112031
$ws.Range("G97").Value = "Poroto verde"
$ws.Range("H97").Value = "Sin especificar"
$ws.Range("I97").Value = "Primera"
$ws.Range("J97").Value = 1330
$ws.Range("K97").Value = 700
$ws.Range("L97").Value = 800
$ws.Range("M97").Value = 751
$ws.Range("N97").Value = "$/kilo"
$ws.Range("O97").Value = "Región de Arica y Parinacota"
$ws.Range("P97").Value = 751
$ws.Range("Q97").Value = 1
$ws.Range("R97").Value = "Hortaliza"
